# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (period) rows for the two workers (DAYANA CAROLINA PEREZ
# HURTADO / CC 1143386199 and ORLANDO RAMON BELTRAN RODRIGUEZ / CC 92070402)
# are re-sequenced so the periods 2404-2410 interleave worker-by-worker in
# ascending period order, instead of each worker's block being listed
# separately in descending period order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16 keeps DAYANA, but now for period 2404 (was 2408)
$ws.Range("E16").Value = "2404"

# Row 17 switches from DAYANA/2407 to ORLANDO/2404
$ws.Range("C17").Value = "92070402"
$ws.Range("D17").Value = "ORLANDO RAMON BELTRAN RODRIGUEZ"
$ws.Range("E17").Value = "2404"

# Row 18 keeps DAYANA, period becomes 2405 (was 2406)
$ws.Range("E18").Value = "2405"

# Row 19 switches from DAYANA to ORLANDO, period stays 2405
$ws.Range("C19").Value = "92070402"
$ws.Range("D19").Value = "ORLANDO RAMON BELTRAN RODRIGUEZ"

# Row 20 keeps DAYANA, period becomes 2406 (was 2404)
$ws.Range("E20").Value = "2406"

# Row 21 keeps ORLANDO, period becomes 2406 (was 2410) and Valor Mora back to 52000
$ws.Range("E21").Value = "2406"
$ws.Range("F21").Value = 52000

# Row 22 switches from ORLANDO/2409 to DAYANA/2407
$ws.Range("C22").Value = "1143386199"
$ws.Range("D22").Value = "DAYANA CAROLINA PEREZ HURTADO"
$ws.Range("E22").Value = "2407"

# Row 23 keeps ORLANDO, period becomes 2407 (was 2408)
$ws.Range("E23").Value = "2407"

# Row 24 switches from ORLANDO/2407 to DAYANA/2408
$ws.Range("C24").Value = "1143386199"
$ws.Range("D24").Value = "DAYANA CAROLINA PEREZ HURTADO"
$ws.Range("E24").Value = "2408"

# Row 25 keeps ORLANDO, period becomes 2408 (was 2406)
$ws.Range("E25").Value = "2408"

# Row 26 keeps ORLANDO, period becomes 2409 (was 2405)
$ws.Range("E26").Value = "2409"

# Row 27 keeps ORLANDO, period becomes 2410 (was 2404), Valor Mora back to 15600
$ws.Range("E27").Value = "2410"
$ws.Range("F27").Value = 15600
